# Subindo CSS e imagens faltantes
# Append the new rows (ids 55-64) of user data to the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(55, "zinter", "meyerrfd44@gmail.com", "66666wdafd"),
    @(56, "srluiz", "luiz224@gmail.com", "luiz3747343"),
    @(57, "Pedro Alves cabral", "cabral@gmail", "cabral"),
    @(58, "Maria ", "Ju@pet", "1212"),
    @(59, "erere", "mekdafeafeaanics153@gmail.com", "erfe"),
    @(60, "JSONvc", "mferreiradejesus11@gmail.com", "mbnv bcxvc"),
    @(61, "Gustavo", "gustavo.vannucchi.ungari@gmail.com", "Narcelo2026="),
    @(62, "Sroliver", "sroliver555@gmail.com", "sroliver"),
    @(63, "tkinterr", "102030@gmail.com", "102030"),
    @(64, "zcrustess", "mferredsdaejesus11@gmail.com", "mnbv")
)

$startRow = 56
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Column D ("senha") must stay text even when the password looks like
    # a plain number (e.g. "1212", "102030"), matching the source data
    # where every value in this column is stored as a string. A leading
    # apostrophe forces Excel to keep such values as text instead of
    # auto-converting them to numbers.
    $senha = $row[3]
    if ($senha -match '^[0-9]+$') {
        $ws.Cells.Item($r, 4).Value = "'" + $senha
    } else {
        $ws.Cells.Item($r, 4).Value = $senha
    }
}
